# Update "合肥-漫展信息" workbook to the next scrape snapshot.
#
# Sheet 1 "展览" (Exhibitions): "想去人数" (F column) counters ticked up/down
#   for several still-open exhibitions.
# Sheet 2 "演出" (Performances): the 2024-08-02 Luke Thompson gig has
#   finished/expired and is dropped from the list; the remaining two rows
#   shift up.
# Sheet 3 "本地生活" (Local life): no changes.
# Sheet 4 "全部类型" (All types) is the same union, minus the same expired
#   Luke Thompson row, with the same F-column ticks as sheet 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 - 展览: bump "想去人数" (F column) counters.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 7641
$ws1.Range("F3").Value = 293
$ws1.Range("F4").Value = 28
$ws1.Range("F6").Value = 4304
$ws1.Range("F7").Value = 327
$ws1.Range("F8").Value = 600
$ws1.Range("F10").Value = 679
$ws1.Range("F11").Value = 160

# ---------------------------------------------------------------------
# Sheet 2 - 演出: the 2024-08-02 Luke Thompson show drops off the list.
# Deleting row 2 shifts rows 3/4 up into rows 2/3.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(2).Delete()

# Column A is a plain positional row counter (1, 2, 3, ...), independent of
# the row's content - renumber it after the delete/shift.
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2

# ---------------------------------------------------------------------
# Sheet 4 - 全部类型: same Luke Thompson row disappears; everything below
# shifts up one row.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(2).Delete()

# Renumber the positional row counter in column A back to 1..12.
for ($i = 0; $i -lt 12; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $i + 1
}

# Re-apply the same "想去人数" (F column) ticks as sheet 1, now that the
# rows have shifted up by one.
$ws4.Range("F2").Value = 7641
$ws4.Range("F4").Value = 293
$ws4.Range("F5").Value = 28
$ws4.Range("F7").Value = 4304
$ws4.Range("F8").Value = 327
$ws4.Range("F9").Value = 600
$ws4.Range("F11").Value = 679
$ws4.Range("F13").Value = 160
